# Insert a new weekly record (Primera/Segunda pair) for "Betarraga" at the top
# of the dated data block (row 134), pushing all subsequent rows down by two
# rows. The new pair repeats the most-recent price figures with a new date
# (2022-01-25 == Excel serial 44586).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 134: down-ward by inserting two fresh rows above the current
# row 134 (this carries the existing formatting of row 134 onto the new
# rows, same as Excel does for a normal "Insert Copied Cells"/"Insert Rows").
$ws.Rows.Item(134).Resize(2).Insert()

# New row 134 - "Primera"
$ws.Cells.Item(134, 1).Value = 8
$ws.Cells.Item(134, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(134, 3).Value = "Coquimbo"
$ws.Cells.Item(134, 4).Value = 44586
$ws.Cells.Item(134, 5).Value = 4
$ws.Cells.Item(134, 6).Value = 100114014
$ws.Cells.Item(134, 7).Value = "Betarraga"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 3100
$ws.Cells.Item(134, 11).Value = 450
$ws.Cells.Item(134, 12).Value = 500
$ws.Cells.Item(134, 13).Value = 475
$ws.Cells.Item(134, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(134, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(134, 16).Value = 158
$ws.Cells.Item(134, 17).Value = 3
$ws.Cells.Item(134, 18).Value = "Hortaliza"

# New row 135 - "Segunda"
$ws.Cells.Item(135, 1).Value = 8
$ws.Cells.Item(135, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(135, 3).Value = "Coquimbo"
$ws.Cells.Item(135, 4).Value = 44586
$ws.Cells.Item(135, 5).Value = 4
$ws.Cells.Item(135, 6).Value = 100114014
$ws.Cells.Item(135, 7).Value = "Betarraga"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Segunda"
$ws.Cells.Item(135, 10).Value = 1520
$ws.Cells.Item(135, 11).Value = 350
$ws.Cells.Item(135, 12).Value = 400
$ws.Cells.Item(135, 13).Value = 375
$ws.Cells.Item(135, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(135, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(135, 16).Value = 125
$ws.Cells.Item(135, 17).Value = 3
$ws.Cells.Item(135, 18).Value = "Hortaliza"
